$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value, applied as literal text
# (NumberFormat "@" forces text storage so numeric-looking strings like
# "1.00" or "5.61" keep their exact textual representation instead of
# being coerced to a floating point number; Style reset afterwards keeps
# the cell formatting identical to the untouched default cells).
$updates = [ordered]@{
    'D2' = '64.379.01'
    'E2' = '  -3.05%  '
    'D3' = '3.176.45'
    'E3' = '  -4.53%  '
    'E4' = '  -0.03%  '
    'D5' = '569.83'
    'E5' = '  -2.80%  '
    'D6' = '169.03'
    'E6' = '  -7.68%  '
    'E7' = '  -5.96%  '
    'E8' = '  -0.02%  '
    'D9' = '3.175.86'
    'E9' = '  -4.54%  '
    'E10' = '  -4.99%  '
    'D11' = '6.76'
    'E11' = '  -0.16%  '
    'E12' = '  -4.38%  '
    'D13' = '3.724.11'
    'E13' = '  -4.60%  '
    'E14' = '  -2.03%  '
    'D15' = '64.428.60'
    'E15' = '  -3.02%  '
    'D16' = '25.35'
    'E16' = '  -3.89%  '
    'E17' = '  -3.71%  '
    'D18' = '3.176.63'
    'E18' = '  -4.20%  '
    'D19' = '419.59'
    'E19' = '  -2.87%  '
    'B20' = 'Chainlink'
    'C20' = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
    'D20' = '12.83'
    'E20' = '  -3.65%  '
    'B21' = 'Polkadot'
    'C21' = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
    'D21' = '5.36'
    'E21' = '  -3.40%  '
    'D22' = '7.04'
    'E22' = '  -5.41%  '
    'D23' = '1.00'
    'E23' = '  -0.19%  '
    'D24' = '70.21'
    'E24' = '  -2.79%  '
    'D25' = '0.203'
    'E25' = '  +2.21%  '
    'E26' = '  -5.80%  '
    'E27' = '  -8.02%  '
    'D28' = '8.73'
    'E28' = '  -3.33%  '
    'D29' = '0.999'
    'E29' = '  -0.06%  '
    'D30' = '1.83'
    'E30' = '  -6.44%  '
    'D31' = '21.76'
    'E31' = '  -3.04%  '
    'E32' = '  -0.12%  '
    'D33' = '5.03'
    'E33' = '  -3.87%  '
    'D34' = '6.32'
    'E34' = '  -4.81%  '
    'E35' = '  -6.02%  '
    'D36' = '157.35'
    'E36' = '  -1.52%  '
    'D37' = '1.36'
    'E37' = '  -7.22%  '
    'D38' = '2.707.74'
    'E38' = '  -6.19%  '
    'D39' = '1.70'
    'E39' = '  -7.33%  '
    'D40' = '24.34'
    'E40' = '  -8.99%  '
    'D41' = '4.16'
    'E41' = '  -4.19%  '
    'D42' = '39.17'
    'E42' = '  -2.96%  '
    'E43' = '  -7.78%  '
    'E44' = '  -7.20%  '
    'D45' = '5.61'
    'E45' = '  -7.22%  '
    'D47' = '21.66'
    'E47' = '  -7.81%  '
    'D48' = '292.28'
    'E48' = '  -8.33%  '
    'E49' = '  -0.04%  '
    'B50' = 'dogwifhat'
    'C50' = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
    'D50' = '2.01'
    'E50' = '  -13.69%  '
    'B51' = 'Stellar'
    'C51' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'D51' = '2.01'
    'E51' = '  -5.94%  '
}

foreach ($ref in $updates.Keys) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $updates[$ref]
    $c.Style = "Normal"
}
